# Reserve Margin.xlsx update: Indonesia -> U.S. NERC data (v3.3.1)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("RM")

# ---------------------------------------------------------------------------
# "About" sheet (sheet1)
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "RM Reserve Margin"
$ws1.Range("A3").Value = "Source:"
$ws1.Range("B3").Value = "North American Electrict Reliability Coorporation "
$ws1.Range("B4").Value = 2015
$ws1.Range("B5").Value = "2015 Summer Reliability Assessment"
$ws1.Range("B6").Value = "http://www.nerc.com/pa/RAPA/ra/Reliability%20Assessments%20DL/2015_Summer_Reliability_Assessment.pdf"
$ws1.Range("B7").Value = "p.3, Table 1: Projected Demand, Resources, and Planning Reserve Margins, NERC Reference Margin Level (%)"
$ws1.Range("A9").Value = "Notes"
$ws1.Range("A10").Value = "The reserve margin (difference between the total generation available and the forecasted peak demand) in the U.S. "
$ws1.Range("A11").Value = "dataset doesn't vary by year, but the RM Reserve Margin variable is a time series to support countries that project "
$ws1.Range("A12").Value = "changes in future reserve margin by year."

# ---------------------------------------------------------------------------
# "RM" sheet (sheet2)
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "(dimensionless)"
$ws2.Range("A2").Value = "Reserve Margin"
$ws2.Range("B2:AK2").Value = 0.1412

# ---------------------------------------------------------------------------
# Selections / active cells to mirror the final saved view state
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A2").Select()

$ws1.Activate()
$ws1.Range("F17").Select()

Write-Output "Reserve Margin workbook updated"
